$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("D2").Value = 44335
    $ws.Range("J2").Value = 1000
    $ws.Range("K2").Value = 12000
    $ws.Range("L2").Value = 13000
    $ws.Range("M2").Value = 12500
    $ws.Range("P2").Value = 500
    # Row 3
    $ws.Range("D3").Value = 44441
    $ws.Range("J3").Value = 1100
    $ws.Range("K3").Value = 11000
    $ws.Range("L3").Value = 12000
    $ws.Range("M3").Value = 11500
    $ws.Range("P3").Value = 460
    # Row 4
    $ws.Range("D4").Value = 44504
    $ws.Range("J4").Value = 700
    $ws.Range("K4").Value = 6000
    $ws.Range("L4").Value = 7000
    $ws.Range("M4").Value = 6500
    $ws.Range("P4").Value = 260
    # Row 5
    $ws.Range("D5").Value = 44455
    $ws.Range("J5").Value = 600
    $ws.Range("K5").Value = 9000
    $ws.Range("L5").Value = 10000
    $ws.Range("M5").Value = 9500
    $ws.Range("P5").Value = 380
    # Row 6
    $ws.Range("D6").Value = 44419
    $ws.Range("J6").Value = 1100
    $ws.Range("K6").Value = 11000
    $ws.Range("L6").Value = 12000
    $ws.Range("M6").Value = 11500
    $ws.Range("P6").Value = 460
    # Row 7
    $ws.Range("D7").Value = 44412
    $ws.Range("J7").Value = 1000
    $ws.Range("K7").Value = 10000
    $ws.Range("L7").Value = 11000
    $ws.Range("M7").Value = 10500
    $ws.Range("P7").Value = 420
    # Row 8
    $ws.Range("D8").Value = 44462
    # Row 9
    $ws.Range("D9").Value = 44503
    $ws.Range("J9").Value = 760
    $ws.Range("K9").Value = 5000
    $ws.Range("L9").Value = 6000
    $ws.Range("M9").Value = 5500
    $ws.Range("P9").Value = 220
    # Row 10
    $ws.Range("D10").Value = 44399
    $ws.Range("K10").Value = 9000
    $ws.Range("L10").Value = 10000
    $ws.Range("M10").Value = 9500
    $ws.Range("P10").Value = 380
    # Row 11
    $ws.Range("D11").Value = 44377
    $ws.Range("K11").Value = 9000
    $ws.Range("L11").Value = 10000
    $ws.Range("M11").Value = 9500
    $ws.Range("P11").Value = 380
    # Row 12
    $ws.Range("D12").Value = 44435
    $ws.Range("J12").Value = 600
    $ws.Range("K12").Value = 10000
    $ws.Range("L12").Value = 11000
    $ws.Range("M12").Value = 10500
    $ws.Range("P12").Value = 420
    # Row 13
    $ws.Range("D13").Value = 44356
    $ws.Range("J13").Value = 1000
    # Row 14
    $ws.Range("D14").Value = 44448
    $ws.Range("J14").Value = 800
    $ws.Range("K14").Value = 10000
    $ws.Range("L14").Value = 12000
    $ws.Range("M14").Value = 11000
    $ws.Range("P14").Value = 440
    # Row 15
    $ws.Range("D15").Value = 44336
    $ws.Range("J15").Value = 1200
    # Row 16
    $ws.Range("D16").Value = 44343
    $ws.Range("J16").Value = 500
    # Row 17
    $ws.Range("D17").Value = 44363
    $ws.Range("J17").Value = 900
    $ws.Range("K17").Value = 11000
    $ws.Range("L17").Value = 12000
    $ws.Range("M17").Value = 11500
    $ws.Range("P17").Value = 460
    # Row 18
    $ws.Range("D18").Value = 44406
    $ws.Range("J18").Value = 800
    $ws.Range("K18").Value = 10000
    $ws.Range("L18").Value = 11000
    $ws.Range("M18").Value = 10500
    $ws.Range("P18").Value = 420
    # Row 19
    $ws.Range("D19").Value = 44392
    $ws.Range("K19").Value = 9000
    $ws.Range("L19").Value = 10000
    $ws.Range("M19").Value = 9500
    $ws.Range("P19").Value = 380
    # Row 20
    $ws.Range("D20").Value = 44426
    $ws.Range("J20").Value = 500
    $ws.Range("K20").Value = 11000
    $ws.Range("M20").Value = 11500
    $ws.Range("P20").Value = 460
    # Row 21
    $ws.Range("D21").Value = 44384
    $ws.Range("J21").Value = 700
    $ws.Range("K21").Value = 8000
    $ws.Range("L21").Value = 9000
    $ws.Range("M21").Value = 8500
    $ws.Range("P21").Value = 340
    # Row 22
    $ws.Range("D22").Value = 44483
    $ws.Range("J22").Value = 1200
    $ws.Range("K22").Value = 4000
    $ws.Range("L22").Value = 5000
    $ws.Range("M22").Value = 4500
    $ws.Range("P22").Value = 180
    # Row 23
    $ws.Range("D23").Value = 44349
    $ws.Range("J23").Value = 600
    $ws.Range("L23").Value = 12000
    $ws.Range("M23").Value = 11000
    $ws.Range("P23").Value = 440
    # Row 24
    $ws.Range("D24").Value = 44482
    $ws.Range("J24").Value = 1600
    $ws.Range("K24").Value = 4000
    $ws.Range("L24").Value = 5000
    $ws.Range("M24").Value = 4500
    $ws.Range("P24").Value = 180
    # Row 25
    $ws.Range("D25").Value = 44364
    $ws.Range("J25").Value = 700
    # Row 26
    $ws.Range("D26").Value = 44469
    $ws.Range("J26").Value = 600
    $ws.Range("K26").Value = 5000
    $ws.Range("L26").Value = 6000
    $ws.Range("M26").Value = 5500
    $ws.Range("P26").Value = 220
    # Row 27
    $ws.Range("D27").Value = 44427
    $ws.Range("J27").Value = 360
    $ws.Range("K27").Value = 10000
    $ws.Range("L27").Value = 11000
    $ws.Range("M27").Value = 10500
    $ws.Range("P27").Value = 420
    # Row 28
    $ws.Range("D28").Value = 44475
    $ws.Range("J28").Value = 1200
    $ws.Range("K28").Value = 5000
    $ws.Range("L28").Value = 6000
    $ws.Range("M28").Value = 5500
    $ws.Range("P28").Value = 220
    # Row 29
    $ws.Range("D29").Value = 44468
    $ws.Range("J29").Value = 700
    # Row 30
    $ws.Range("D30").Value = 44489
    $ws.Range("J30").Value = 1200
    # Row 31
    $ws.Range("D31").Value = 44476
    $ws.Range("K31").Value = 5000
    $ws.Range("L31").Value = 6000
    $ws.Range("M31").Value = 5500
    $ws.Range("P31").Value = 220
    # Row 32
    $ws.Range("D32").Value = 44447
    $ws.Range("J32").Value = 1000
    # Row 33
    $ws.Range("D33").Value = 44328
    $ws.Range("J33").Value = 900
    $ws.Range("K33").Value = 11000
    $ws.Range("L33").Value = 12000
    $ws.Range("M33").Value = 11500
    $ws.Range("P33").Value = 460
    # Row 34
    $ws.Range("D34").Value = 44434
    $ws.Range("J34").Value = 600
    # Row 35
    $ws.Range("D35").Value = 44490
    $ws.Range("J35").Value = 400
    $ws.Range("K35").Value = 5000
    $ws.Range("L35").Value = 6000
    $ws.Range("M35").Value = 5500
    $ws.Range("P35").Value = 220
    # Row 36
    $ws.Range("D36").Value = 44391
    $ws.Range("J36").Value = 500
    # Row 37
    $ws.Range("D37").Value = 44510
    $ws.Range("J37").Value = 1300
    $ws.Range("K37").Value = 6000
    $ws.Range("L37").Value = 7000
    $ws.Range("M37").Value = 6500
    $ws.Range("P37").Value = 260
    # Row 38
    $ws.Range("D38").Value = 44420
    $ws.Range("J38").Value = 1000
    $ws.Range("K38").Value = 10000
    $ws.Range("L38").Value = 11000
    $ws.Range("M38").Value = 10500
    $ws.Range("P38").Value = 420
    # Row 39
    $ws.Range("D39").Value = 44385
    $ws.Range("J39").Value = 600
    $ws.Range("K39").Value = 8000
    $ws.Range("L39").Value = 9000
    $ws.Range("M39").Value = 8500
    $ws.Range("P39").Value = 340
    # Row 40
    $ws.Range("D40").Value = 44308
    $ws.Range("J40").Value = 400
    $ws.Range("K40").Value = 11000
    $ws.Range("L40").Value = 12000
    $ws.Range("M40").Value = 11500
    $ws.Range("P40").Value = 460
    # Row 41
    $ws.Range("D41").Value = 44413
    $ws.Range("J41").Value = 1200
    $ws.Range("K41").Value = 10000
    $ws.Range("L41").Value = 11000
    $ws.Range("M41").Value = 10500
    $ws.Range("P41").Value = 420
    # Row 42
    $ws.Range("D42").Value = 44329
    $ws.Range("J42").Value = 1000
    $ws.Range("K42").Value = 12000
    $ws.Range("L42").Value = 13000
    $ws.Range("M42").Value = 12500
    $ws.Range("P42").Value = 500
    # Row 43
    $ws.Range("D43").Value = 44398
    $ws.Range("J43").Value = 400
    # Row 44
    $ws.Range("D44").Value = 44461
    $ws.Range("J44").Value = 1100
    $ws.Range("K44").Value = 9000
    $ws.Range("L44").Value = 10000
    $ws.Range("M44").Value = 9500
    $ws.Range("P44").Value = 380
    # Row 45
    $ws.Range("D45").Value = 44371
    $ws.Range("J45").Value = 500
    # Row 46
    $ws.Range("D46").Value = 44454
    $ws.Range("J46").Value = 800
    $ws.Range("K46").Value = 9000
    $ws.Range("L46").Value = 10000
    $ws.Range("M46").Value = 9500
    $ws.Range("P46").Value = 380
